$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1021
$ws.Range("F4").Value = 248
$ws.Range("F5").Value = 18
$ws.Range("F6").Value = 459
$ws.Range("F7").Value = 742
$ws.Range("F8").Value = 256
$ws.Range("F11").Value = 413
$ws.Range("F12").Value = 218
$ws.Range("F13").Value = 84
$ws.Range("F14").Value = 858
$ws.Range("F15").Value = 117
$ws.Range("F16").Value = 2003
$ws.Range("F17").Value = 494
$ws.Range("F18").Value = 7619
$ws.Range("F19").Value = 559

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 3
$ws.Range("F8").Value = 126
$ws.Range("F10").Value = 8

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5543
$ws.Range("F4").Value = 392

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 5543
$ws.Range("F5").Value = 392
$ws.Range("F7").Value = 1021
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = 248
$ws.Range("F11").Value = 18
$ws.Range("F12").Value = 459
$ws.Range("F13").Value = 742
$ws.Range("F14").Value = 256
$ws.Range("F18").Value = 413
$ws.Range("F19").Value = 218
$ws.Range("F21").Value = 84
$ws.Range("F23").Value = 858
$ws.Range("F24").Value = 117
$ws.Range("F25").Value = 126
$ws.Range("F26").Value = 2003
$ws.Range("F27").Value = 494
$ws.Range("F28").Value = 7619
$ws.Range("F30").Value = 8
$ws.Range("F31").Value = 559
